$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (new TPM recalculated values)
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 5.027277643914444
$ws.Range("R2").Value = 45.24549879523001
$ws.Range("S2").Value = 0.04635500474236593
$ws.Range("T2").Value = 0.04635500474236593

# Row 3 updates
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("S3").Value = 0.6912512390256352
$ws.Range("T3").Value = 0.6912512390256351

# Row 4 updates
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 28.45704087270333
$ws.Range("R4").Value = 256.11336785433
$ws.Range("S4").Value = 0.2623937562319988
$ws.Range("T4").Value = 0.2623937562319988
